# Updated scenarios (removed H2), and set initial capacity for Transmission and Distribution
#
# The "EC_H2" scenario row (row 4) in both the lowBio and highBio sheets is
# removed entirely; remaining rows shift up to close the gap.

$wb = $excel.ActiveWorkbook

$lowBio = $wb.Worksheets.Item("lowBio")
$highBio = $wb.Worksheets.Item("highBio")

# Remove the EC_H2 scenario row on the highBio sheet first...
$highBio.Activate()
$highBio.Rows.Item(4).EntireRow.Delete()
$highBio.Range("A4:XFD4").Select()

# ...then on the lowBio sheet, leaving lowBio as the active/selected sheet
# (matches the original workbook's active tab).
$lowBio.Activate()
$lowBio.Rows.Item(4).EntireRow.Delete()

# Leave the selection where the row used to be, same as the row that slid
# up into its place after the delete.
$lowBio.Range("A4:XFD4").Select()
